$wb = $excel.ActiveWorkbook

$dateFormat = "YYYY-MM-DD HH:MM:SS"

function Set-LogRow {
    param(
        $ws,
        $row,
        $a,
        $b,
        $c,
        $d,
        $e,
        $f,
        $g,
        $h,
        $i
    )

    $ws.Cells.Item($row, 1).Value = $a
    $ws.Cells.Item($row, 1).NumberFormat = $dateFormat
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = [double]$g
    $ws.Cells.Item($row, 8).Value = $h
    $ws.Cells.Item($row, 9).Value = $i
}

# ---- Sheet 1: DE_LFT_#1 ----
$ws1 = $wb.Worksheets.Item("DE_LFT_#1")
Set-LogRow $ws1 159 45945.43409722222 "0x01,0x7c" "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0," "0x00,0xE4" "0x14" 380 "7.598631275147109e+23" 228 14
Set-LogRow $ws1 160 45946.43524305556 "0x01,0x7c" "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0," "0x00,0xE4" "0x14" 380 "7.598631275147109e+23" 228 14

# ---- Sheet 2: DE_LFT_#2 ----
$ws2 = $wb.Worksheets.Item("DE_LFT_#2")
Set-LogRow $ws2 159 45945.43409722222 "0x01,0x7c" "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78," "0x00,0xEC" "0xe" 380 "5.68432987514711e+23" 236 14
Set-LogRow $ws2 160 45946.43524305556 "0x01,0x7c" "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78," "0x00,0xEC" "0xe" 380 "5.68432987514711e+23" 236 14

# ---- Sheet 3: DE_PLT_#1 ----
$ws3 = $wb.Worksheets.Item("DE_PLT_#1")
Set-LogRow $ws3 159 45945.43409722222 "0x00,0x82" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c," "0x00,0x6A" "0x7" 130 "5.68631262647114e+23" 106 7
Set-LogRow $ws3 160 45946.43524305556 "0x00,0x82" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c," "0x00,0x6A" "0x7" 130 "5.68631262647114e+23" 106 7

# ---- Sheet 4: DE_PLT_#2 ----
$ws4 = $wb.Worksheets.Item("DE_PLT_#2")
Set-LogRow $ws4 159 45945.43409722222 "0x00,0x82" "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c," "0x00,0x68" "0x3" 130 "9.85046333984776e+23" 104 3
Set-LogRow $ws4 160 45946.43524305556 "0x00,0x82" "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c," "0x00,0x68" "0x3" 130 "9.85046333984776e+23" 104 3
